$d = $word.ActiveDocument
$wordMl = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# Paragraph 1: "Konspekt lekcji 4" -- drop the strike-through that was
# carried on the paragraph mark only (the visible text itself never had
# strike-through applied to it).
$p1 = $d.Paragraphs.Item(1)
$xml1 = '<w:p xmlns:w="' + $wordMl + '"><w:r><w:t>Konspekt lekcji 4</w:t></w:r></w:p>'
$null = $p1.Range.InsertXML($xml1)

# Paragraph 2: used to hold only the hidden _GoBack bookmark. Give it the
# new text run "Funkcja i rozwiazania do lekcji 4" (placed before the
# bookmark, exactly where it was typed) and move the strike-through onto
# this paragraph's own mark instead.
$p2 = $d.Paragraphs.Item(2)
$text2 = "Funkcja i rozwi" + [char]0x0105 + "zania do lekcji 4"
$xml2 = '<w:p xmlns:w="' + $wordMl + '">' +
        '<w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>' +
        '<w:r><w:t>' + $text2 + '</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '</w:p>'
$null = $p2.Range.InsertXML($xml2)

# Add a brand-new, completely empty paragraph right after paragraph 2
# (before the final section break).
$p2 = $d.Paragraphs.Item(2)
$null = $p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$xml3 = '<w:p xmlns:w="' + $wordMl + '"/>'
$null = $p3.Range.InsertXML($xml3)

$d.Save()
Write-Output "OK"
